$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# Sheet "展览" (sheet1) updates to column F
$ws1.Range("F3").Value = 83
$ws1.Range("F4").Value = 1519
$ws1.Range("F6").Value = 1075
$ws1.Range("F7").Value = 11144
$ws1.Range("F8").Value = 86
$ws1.Range("F9").Value = 33
$ws1.Range("F10").Value = 325
$ws1.Range("F11").Value = 1071
$ws1.Range("F12").Value = 762
$ws1.Range("F13").Value = 12253
$ws1.Range("F14").Value = 12831
$ws1.Range("F16").Value = 127
$ws1.Range("F21").Value = 51

# Sheet "全部类型" (sheet4) updates to column F
$ws4.Range("F4").Value = 83
$ws4.Range("F5").Value = 1519
$ws4.Range("F7").Value = 1075
$ws4.Range("F8").Value = 11144
$ws4.Range("F9").Value = 86
$ws4.Range("F10").Value = 33
$ws4.Range("F11").Value = 325
$ws4.Range("F12").Value = 1071
$ws4.Range("F13").Value = 762
$ws4.Range("F14").Value = 12253
$ws4.Range("F15").Value = 12831
$ws4.Range("F17").Value = 127
$ws4.Range("F22").Value = 51
